$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1177.625
$ws.Range("I19").Value = 1232.2
$ws.Range("J19").Value = 1086.6666
$ws.Range("K19").Value = 1232.2
$ws.Range("L19").Value = 1086.6666
$ws.Range("M19").Value = -1057.2
$ws.Range("N19").Value = -1436.6666
$ws.Range("H33").Value = 1397.762
$ws.Range("I33").Value = 1318.4736
$ws.Range("K33").Value = 1318.4736
$ws.Range("M33").Value = -1089.4736
$ws.Range("H62").Value = 381097.4
$ws.Range("I62").Value = 381097.4
$ws.Range("K62").Value = 381097.4
$ws.Range("M62").Value = -380473.4
$ws.Range("H65").Value = 381097.4
$ws.Range("I65").Value = 381097.4
$ws.Range("K65").Value = 1905487
$ws.Range("M65").Value = -1902367
$ws.Range("H86").Value = 1534.3334
$ws.Range("I86").Value = 1381.2
$ws.Range("K86").Value = 1381.2
$ws.Range("M86").Value = -258.2
$ws.Range("H89").Value = 1534.3334
$ws.Range("I89").Value = 1381.2
$ws.Range("K89").Value = 6906
$ws.Range("M89").Value = -1290
$ws.Range("H106").Value = 118488.555
$ws.Range("I106").Value = 148056.72
$ws.Range("J106").Value = 15000
$ws.Range("K106").Value = 148056.72
$ws.Range("L106").Value = 15000
$ws.Range("M106").Value = -147425.72
$ws.Range("N106").Value = -16262
$ws.Range("H107").Value = 1579.4
$ws.Range("I107").Value = 1666.3334
$ws.Range("J107").Value = 1449
$ws.Range("K107").Value = 1666.3334
$ws.Range("L107").Value = 1449
$ws.Range("M107").Value = 253.6666
$ws.Range("N107").Value = -5289
$ws.Range("H113").Value = 5668
$ws.Range("I113").Value = 1999.6666
$ws.Range("J113").Value = 7502.1665
$ws.Range("K113").Value = 1999.6666
$ws.Range("L113").Value = 7502.1665
$ws.Range("M113").Value = 1254.3334
$ws.Range("N113").Value = -14010.1665
$ws.Range("H137").Value = 2341.0356
$ws.Range("I137").Value = 2352.125
$ws.Range("K137").Value = 7056.375
$ws.Range("M137").Value = -4506.375
$ws.Range("H138").Value = 1839574
$ws.Range("I138").Value = 1104.9231
$ws.Range("J138").Value = 2835411.5
$ws.Range("K138").Value = 3314.7693
$ws.Range("L138").Value = 8506234.5
$ws.Range("M138").Value = 1825.2307
$ws.Range("N138").Value = -8516514.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2482.68
$ws.Range("I32").Value = 2511.7568
$ws.Range("K32").Value = 2511.7568
$ws.Range("M32").Value = -2224.7568
$ws.Range("H45").Value = 5008.857
$ws.Range("I45").Value = 1682
$ws.Range("K45").Value = 1682
$ws.Range("M45").Value = -1305
$ws.Range("H110").Value = 2395.4167
$ws.Range("I110").Value = 1718.375
$ws.Range("K110").Value = 1718.375
$ws.Range("M110").Value = 326.625
$ws.Range("H122").Value = 2502.077
$ws.Range("I122").Value = 2377.8
$ws.Range("K122").Value = 7133.400000000001
$ws.Range("M122").Value = -4683.400000000001
$ws.Range("H132").Value = 4006.2163
$ws.Range("I132").Value = 3046
$ws.Range("J132").Value = 8121.4287
$ws.Range("K132").Value = 9138
$ws.Range("L132").Value = 24364.2861
$ws.Range("M132").Value = -6608
$ws.Range("N132").Value = -29424.2861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 2012.8
$ws.Range("J64").Value = 2021.3334
$ws.Range("L64").Value = 2021.3334
$ws.Range("N64").Value = -2471.3334
$ws.Range("H67").Value = 2012.8
$ws.Range("J67").Value = 2021.3334
$ws.Range("L67").Value = 2021.3334
$ws.Range("N67").Value = -3581.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 364.75
$ws.Range("J7").Value = 201
$ws.Range("L7").Value = 201
$ws.Range("N7").Value = -427
$ws.Range("H31").Value = 4760.294
$ws.Range("I31").Value = 3999
$ws.Range("J31").Value = 4994.5386
$ws.Range("K31").Value = 3999
$ws.Range("L31").Value = 4994.5386
$ws.Range("M31").Value = -3704
$ws.Range("N31").Value = -5584.5386
$ws.Range("H34").Value = 4760.294
$ws.Range("I34").Value = 3999
$ws.Range("J34").Value = 4994.5386
$ws.Range("K34").Value = 3999
$ws.Range("L34").Value = 4994.5386
$ws.Range("M34").Value = -3797
$ws.Range("N34").Value = -5398.5386
$ws.Range("H41").Value = 984.9
$ws.Range("I41").Value = 984.9
$ws.Range("K41").Value = 984.9
$ws.Range("M41").Value = -556.9
$ws.Range("H62").Value = 4501.25
$ws.Range("I62").Value = 4501.25
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4501.25
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3877.25
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 4501.25
$ws.Range("I65").Value = 4501.25
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 22506.25
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -19386.25
$ws.Range("N65").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 237.16667
$ws.Range("I2").Value = 266
$ws.Range("J2").Value = 222.75
$ws.Range("K2").Value = 1596
$ws.Range("L2").Value = 1336.5
$ws.Range("M2").Value = -1483
$ws.Range("N2").Value = -1562.5
$ws.Range("H69").Value = 9953.5
$ws.Range("I69").Value = 9953.5
$ws.Range("K69").Value = 29860.5
$ws.Range("M69").Value = -29049.5
$ws.Range("H72").Value = 9953.5
$ws.Range("I72").Value = 9953.5
$ws.Range("K72").Value = 89581.5
$ws.Range("M72").Value = -85525.5
$ws.Range("H97").Value = 2422.3333
$ws.Range("J97").Value = 943.75
$ws.Range("L97").Value = 2831.25
$ws.Range("N97").Value = -3823.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 139.83333
$ws.Range("I2").Value = 162
$ws.Range("J2").Value = 29
$ws.Range("K2").Value = 162
$ws.Range("L2").Value = 29
$ws.Range("M2").Value = -49
$ws.Range("N2").Value = -255
$ws.Range("H15").Value = 96420
$ws.Range("J15").Value = 96420
$ws.Range("L15").Value = 96420
$ws.Range("N15").Value = -96996
$ws.Range("H43").Value = 86817.09
$ws.Range("J43").Value = 99998
$ws.Range("L43").Value = 99998
$ws.Range("N43").Value = -100300
$ws.Range("H81").Value = 96420
$ws.Range("J81").Value = 96420
$ws.Range("L81").Value = 96420
$ws.Range("N81").Value = -98416
$ws.Range("H84").Value = 96420
$ws.Range("J84").Value = 96420
$ws.Range("L84").Value = 289260
$ws.Range("N84").Value = -299244
$ws.Range("H113").Value = 1992.5883
$ws.Range("I113").Value = 1583.5
$ws.Range("K113").Value = 1583.5
$ws.Range("M113").Value = 586.5
$ws.Range("H122").Value = 66669870
$ws.Range("I122").Value = 76926390
$ws.Range("K122").Value = 230779170
$ws.Range("M122").Value = -230776720
$ws.Range("H123").Value = 53153.855
$ws.Range("J123").Value = 53153.855
$ws.Range("L123").Value = 53153.855
$ws.Range("N123").Value = -58053.855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3782.0454
$ws.Range("I22").Value = 756.55554
$ws.Range("J22").Value = 5876.615
$ws.Range("K22").Value = 756.55554
$ws.Range("L22").Value = 5876.615
$ws.Range("M22").Value = -461.55554
$ws.Range("N22").Value = -6466.615
$ws.Range("H27").Value = 3782.0454
$ws.Range("I27").Value = 756.55554
$ws.Range("J27").Value = 5876.615
$ws.Range("K27").Value = 756.55554
$ws.Range("L27").Value = 5876.615
$ws.Range("M27").Value = -649.55554
$ws.Range("N27").Value = -6090.615
$ws.Range("H46").Value = 3928.8845
$ws.Range("J46").Value = 5883
$ws.Range("L46").Value = 5883
$ws.Range("N46").Value = -6259
$ws.Range("H68").Value = 3966.16
$ws.Range("I68").Value = 3939.35
$ws.Range("K68").Value = 3939.35
$ws.Range("M68").Value = -3190.35
$ws.Range("H71").Value = 3966.16
$ws.Range("I71").Value = 3939.35
$ws.Range("K71").Value = 19696.75
$ws.Range("M71").Value = -15952.75
$ws.Range("H122").Value = 5938.091
$ws.Range("I122").Value = 4941.778
$ws.Range("K122").Value = 14825.334
$ws.Range("M122").Value = -12375.334
$ws.Range("H132").Value = 9885.071
$ws.Range("I132").Value = 9884.799999999999
$ws.Range("J132").Value = 9885.75
$ws.Range("K132").Value = 29654.4
$ws.Range("L132").Value = 29657.25
$ws.Range("M132").Value = -27124.4
$ws.Range("N132").Value = -34717.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 6000
$ws.Range("I28").Value = 6000
$ws.Range("K28").Value = 6000
$ws.Range("M28").Value = -5652
$ws.Range("H107").Value = 681.1667
$ws.Range("J107").Value = 785.55554
$ws.Range("L107").Value = 2356.66662
$ws.Range("N107").Value = -6196.66662
$ws.Range("H122").Value = 3162.5
$ws.Range("I122").Value = 3370.9524
$ws.Range("K122").Value = 10112.8572
$ws.Range("M122").Value = -7662.8572
$ws.Range("H126").Value = 6078.75
$ws.Range("J126").Value = 4948.8335
$ws.Range("L126").Value = 14846.5005
$ws.Range("N126").Value = -19786.5005
$ws.Range("H136").Value = 6032.96
$ws.Range("I136").Value = 5269.3887
$ws.Range("K136").Value = 15808.1661
$ws.Range("M136").Value = -13258.1661

Write-Output "Applied 237 cell edits across 8 sheets"